$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.456.17"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "3.520.20"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.85"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.14"
$ws.Range("E6").Value = "  -3.96%  "

$ws.Range("D7").Value = "3.519.34"
$ws.Range("E7").Value = "  -2.29%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  +5.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.70"
$ws.Range("E10").Value = "  -4.77%  "

$ws.Range("E11").Value = "  -4.46%  "

$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").Value = "4.118.99"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("E14").Value = "  -7.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.55"
$ws.Range("E15").Value = "  -4.17%  "

$ws.Range("D16").Value = "3.520.58"
$ws.Range("E16").Value = "  -2.48%  "

$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").Value = "66.331.80"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.84"
$ws.Range("E19").Value = "  -6.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  -3.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.56"
$ws.Range("E21").Value = "  -3.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.05"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.587"
$ws.Range("E23").Value = "  -5.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.16"
$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("D25").Value = "3.662.73"
$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -7.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.88"
$ws.Range("E28").Value = "  -5.40%  "

$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.90"
$ws.Range("E30").Value = "  -5.83%  "

$ws.Range("D32").Value = "3.527.54"
$ws.Range("E32").Value = "  -2.04%  "

$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("E34").Value = "  -5.30%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  -9.57%  "

$ws.Range("E37").Value = "  -4.64%  "

$ws.Range("E38").Value = "  -4.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.30"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.19"
$ws.Range("E40").Value = "  -8.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0816"
$ws.Range("E41").Value = "  -4.76%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.858"
$ws.Range("E42").Value = "  -4.69%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.96"
$ws.Range("E43").Value = "  -5.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.35"
$ws.Range("E44").Value = "  -1.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  -8.11%  "

$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  -8.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.05"
$ws.Range("E48").Value = "  -2.16%  "

$ws.Range("E49").Value = "  -6.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.64"
$ws.Range("E50").Value = "  -6.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.904"
$ws.Range("E51").Value = "  -5.47%  "
